$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.413.30'
$ws.Range('E2').Value = '  -0.95%  '

$ws.Range('D3').Value = '1.562.48'
$ws.Range('E3').Value = '  -1.42%  '

$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.84%  '

$ws.Range('E6').Value = '  -0.92%  '

$ws.Range('E7').Value = '  -0.17%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '21.98'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.13%  '

$ws.Range('E9').Value = '  -1.91%  '

$ws.Range('E10').Value = '  -0.11%  '

$ws.Range('E11').Value = '  -0.10%  '

$ws.Range('D12').Value = '1.785.84'
$ws.Range('E12').Value = '  -1.33%  '

$ws.Range('D13').Value = '1.564.82'
$ws.Range('E13').Value = '  -1.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.81'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.01%  '

$ws.Range('E15').Value = '  -2.71%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.46'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.35%  '

$ws.Range('D17').Value = '27.416.40'
$ws.Range('E17').Value = '  -0.84%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '212.35'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.52%  '

$ws.Range('D19').Value = '0.0₃0688'

$ws.Range('E20').Value = '  -1.02%  '

$ws.Range('E21').Value = '  -0.21%  '

$ws.Range('E22').Value = '  -1.07%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '

$ws.Range('E24').Value = '  +2.87%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.79'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.72%  '

$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.21%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '14.94'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.12%  '

$ws.Range('E29').Value = '  -1.86%  '

$ws.Range('E30').Value = '  -0.16%  '

$ws.Range('E31').Value = '  +0.90%  '

$ws.Range('D33').Value = '1.370.19'
$ws.Range('E33').Value = '  -1.19%  '

$ws.Range('E34').Value = '  +0.44%  '

$ws.Range('E35').Value = '  +0.68%  '

$ws.Range('E36').Value = '  -0.33%  '

$ws.Range('E37').Value = '  -0.89%  '

$ws.Range('E38').Value = '  +1.04%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.529'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.819'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.08%  '

$ws.Range('E41').Value = '  -0.16%  '

$ws.Range('E42').Value = '  -0.33%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.88%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.26%  '

$ws.Range('E45').Value = '  -0.68%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.26%  '

$ws.Range('B47').Value = 'PAXGold'
$ws.Range('C47').Value = 'https://coinranking.com/coin/YRTkUcMi+paxgold-paxg'
$ws.Range('D47').Value = '1.865.97'
$ws.Range('E47').Value = '  -0.25%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.697.97'
$ws.Range('E48').Value = '  -1.37%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.00%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₇0986'
$ws.Range('E50').Value = '  -1.83%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0956'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.75%  '
